# Auto-generated edit script applying the cryptos.xlsx price/symbol update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for numeric-looking text values we must force the cell
# to Text format before assignment (otherwise Excel auto-detects a number),
# then clear the temporary number-format so no stray style is left behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "236.03"
Set-TextValue "D3" "22.41"
Set-TextValue "D4" "5.384"
Set-TextValue "D5" "0.05642"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D6" "6.485"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "3.373"
$ws.Range("E7").Value = "6GateTokenGT"
Set-TextValue "D8" "1.066"
Set-TextValue "D9" "0.7845"
Set-TextValue "D10" "0.1399"
Set-TextValue "D11" "0.07337"
Set-TextValue "D12" "0.03188"
Set-TextValue "D13" "0.02934"
Set-TextValue "D15" "0.001661"
Set-TextValue "D16" "3.254"
Set-TextValue "D17" "0.04764"
Set-TextValue "D19" "0.006226"
Set-TextValue "D20" "0.005106"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "3.862"
Set-TextValue "D25" "0.3323"
Set-TextValue "D26" "0.1055"
Set-TextValue "D40" "0.04104"
Set-TextValue "D41" "0.006989"
Set-TextValue "D42" "0.1038"
Set-TextValue "D43" "0.003263"
Set-TextValue "D44" "0.009938"
Set-TextValue "D45" "0.00005430"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.6753"
Set-TextValue "D48" "0.03901"
Set-TextValue "D49" "0.00002103"
Set-TextValue "D50" "0.01012"
